$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: bump the "Date" value ---------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- 2. Elements sheet: swap the two "Mapping" columns (AK <-> AL) ------
$elements = $wb.Worksheets.Item("Elements")

# Column widths: AK (37) and AL (38) widths are swapped. The stored
# values 24.98046875 / 75.0078125 are Excel's character-width encoding of
# the underlying pixel widths; feed values from the middle of the pixel
# bucket that rounds back to those widths.
$elements.Columns.Item(37).ColumnWidth = 74.15
$elements.Columns.Item(38).ColumnWidth = 24.15

# Swap the cell contents (including the header row) between column AK
# (37) and column AL (38) for every used row.
for ($r = 1; $r -le 17; $r++) {
    $akValue = $elements.Cells.Item($r, 37).Value()
    $alValue = $elements.Cells.Item($r, 38).Value()
    $elements.Cells.Item($r, 37).Value = $alValue
    $elements.Cells.Item($r, 38).Value = $akValue
}
